# ===========================================================================
# Change 1: collapse the split " as " / "the" / " example." runs (with their
# spellStart/gramStart/gramEnd/spellEnd proofErr markers) into a single run by
# replacing the phrase with itself via Find/Replace - Word normalises the run
# structure to one run for the replacement text.
# ===========================================================================
$d = $word.ActiveDocument
[void]$d.Content.Find.Execute(
    " as the example.", $true, $false, $false, $false, $false,
    $true, 1, $false, " as the example.", 2
)

# ===========================================================================
# Change 2: append the new "adventure" paragraph, a blank paragraph, the new
# "Train the Model on New Data" heading and its three mlflow command lines
# right after the existing  mlflow run . -P steps="test_regression_model"  line.
# ===========================================================================
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Contains('test_regression_model')) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "could not locate the test_regression_model anchor paragraph"
}

# Seed 6 brand-new, unstyled paragraphs immediately ahead of the (already
# present) blank paragraph that follows the anchor, so none of them inherit
# the Heading5/Heading5Char formatting used by the anchor paragraph itself.
$seedParagraph = $d.Paragraphs.Item($anchorIndex + 1)
$seedRange = $seedParagraph.Range
for ($i = 0; $i -lt 6; $i++) {
    [void]$seedRange.InsertParagraphBefore()
}

# New narrative paragraph describing the conda.yml / os.path.join / config.yaml changes.
$target = $d.Paragraphs.Item(($anchorIndex + 1))
$targetXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">I had a little bit of an adventure getting this working. I had to change the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>conda.yml</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file in the components/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>test_regression_model</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to remove them installing a git environment. I do not quite understand that. That may bite me. I just adjusted the </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>os.path</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>.join</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to get me in the right directory and updated the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>config.yaml</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file  to reflect the realities of what I had named things on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wandb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target.Range.InsertXML($targetXml)

# Blank spacer paragraph.
$target = $d.Paragraphs.Item(($anchorIndex + 1 + 1))
$targetXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target.Range.InsertXML($targetXml)

# "Train the Model on New Data" - Heading4.
$target = $d.Paragraphs.Item(($anchorIndex + 1 + 2))
$targetXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>Train the Model on New Data</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target.Range.InsertXML($targetXml)

# First mlflow command line - Heading5.
$target = $d.Paragraphs.Item(($anchorIndex + 1 + 3))
$targetXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading5"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>mlflow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> run https://github.com/LindsayMoir/NYC_Predict_Rental_Prices.git \</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target.Range.InsertXML($targetXml)

# Second mlflow command line (-v "1.0.0") - Heading5.
$target = $d.Paragraphs.Item(($anchorIndex + 1 + 4))
$targetXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading5"/></w:pPr><w:r><w:t xml:space="preserve">             -v "1.0.0" \</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target.Range.InsertXML($targetXml)

# Third mlflow command line (-P hydra_options=...) - Heading5.
$target = $d.Paragraphs.Item(($anchorIndex + 1 + 5))
$targetXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading5"/></w:pPr><w:r><w:t xml:space="preserve">             -P </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hydra_options</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>="etl.sample2=''sample2.csv''"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target.Range.InsertXML($targetXml)

Write-Output "paragraph count: $($d.Paragraphs.Count)"
